$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 902, shifting the existing
# rows 902:921 down to 906:925 (matches the growth of the sheet's
# dimension from A1:R921 to A1:R925).
$ws.Rows("902:905").Insert()

# Fill the 4 newly inserted rows with the new weekly price entries.
$newRows = @(
    @{ Row = 902; A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 45239; E = 16; F = 100112020; G = "Tomate"; H = "Larga vida"; I = "Extra";   J = 150; K = 18000; L = 18000; M = 18000; N = "`$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 1000; Q = 18; R = "Hortaliza" },
    @{ Row = 903; A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 45239; E = 16; F = 100112020; G = "Tomate"; H = "Larga vida"; I = "Primera"; J = 300; K = 15000; L = 16000; M = 15500; N = "`$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 861;  Q = 18; R = "Hortaliza" },
    @{ Row = 904; A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 45239; E = 16; F = 100112020; G = "Tomate"; H = "Larga vida"; I = "Primera"; J = 400; K = 8000;  L = 8000;  M = 8000;  N = "`$/caja 10 kilos";    O = "Región de Arica y Parinacota"; P = 800;  Q = 10; R = "Hortaliza" },
    @{ Row = 905; A = 7; B = "Terminal Hortofrutícola Agro Chillán"; C = "Ñuble"; D = 45239; E = 16; F = 100112020; G = "Tomate"; H = "Larga vida"; I = "Segunda"; J = 300; K = 12500; L = 12500; M = 12500; N = "`$/bandeja 18 kilos"; O = "Región de Arica y Parinacota"; P = 694;  Q = 18; R = "Hortaliza" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
}
